$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove "See the [vignette hyperlink] for more details." so the
#    paragraph ends with "...that may benefit. "
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "may benefit. See the vignette for more details.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "may benefit. ", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. "> library(AzureRMR)" console echo - merge "> " and "library(" into a
#    single run (drops the stray grammar-start proofing mark).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "> library(",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "> library(", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3. "> library(SAR)" console echo - merge the whole line into one run
#    (drops the stray grammar-start/grammar-end proofing marks).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "> library(SAR)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "> library(SAR)", 2) | Out-Null

# ---------------------------------------------------------------------------
# 4. "do_operation" - merge "do" and "_operation" into a single run (drops
#    the stray grammar-end proofing mark between them).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "do_operation",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "do_operation", 2) | Out-Null
